$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (shared strings) ---

# A8: "Volume 32   Number  19" -> "...20" (issue number increment)
$cA8 = $ws.Range("A8")
$fullA8 = $cA8.Value()
$idxA8 = $fullA8.IndexOf("19")
$cA8.Characters($idxA8 + 1, 2).Text = "20"

# C9: "Report Covering the Week  5/5/2025  Through  5/11/2025"
#     -> "...5/12/2025  Through  5/18/2025" (next week's reporting period)
$cC9 = $ws.Range("C9")
$fullC9a = $cC9.Value()
$idxC9a = $fullC9a.IndexOf("5/5/2025")
$cC9.Characters($idxC9a + 1, 8).Text = "5/12/2025"
$fullC9b = $cC9.Value()
$idxC9b = $fullC9b.IndexOf("5/11/2025")
$cC9.Characters($idxC9b + 1, 9).Text = "5/18/2025"

# --- Update data table (rows 14-31) with new weekly crime statistics ---

# Cells that change from a blank placeholder to a numeric value (or vice versa)
# need their number format cloned from a same-styled neighboring cell so the
# resulting style matches the rest of the column.

$ws.Range("F22").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1

$ws.Range("M22").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0

$ws.Range("F22").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 1

$ws.Range("F22").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1

$ws.Range("M22").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 0

$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1

$ws.Range("M22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100

$ws.Range("F22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1

$ws.Range("M22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = 0

$ws.Range("F22").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1

$ws.Range("M22").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0

$ws.Range("F22").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 3

$ws.Range("M22").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -33.333333333333

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("D31").Copy()
$ws.Range("C31").PasteSpecial(-4122)

# Plain numeric updates (style unchanged)
$ws.Range("L14").Value = -66.666666666666
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -36.363636363636
$ws.Range("L15").Value = 16.666666666666
$ws.Range("N15").Value = -46.153846153846
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = -20
$ws.Range("M16").Value = -48.387096774193
$ws.Range("N16").Value = -88.965517241379
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -25.714285714285
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 133
$ws.Range("K17").Value = -17.293233082706
$ws.Range("L17").Value = -21.985815602836
$ws.Range("M17").Value = 37.5
$ws.Range("N17").Value = -14.0625
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 41
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = 10.810810810810
$ws.Range("L18").Value = -32.786885245901
$ws.Range("M18").Value = -59.803921568627
$ws.Range("N18").Value = -91.816367265469
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 80
$ws.Range("F19").Value = 20
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 113
$ws.Range("J19").Value = 131
$ws.Range("K19").Value = -13.740458015267
$ws.Range("L19").Value = -0.877192982456
$ws.Range("M19").Value = -6.611570247933
$ws.Range("N19").Value = -42.929292929292
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -37.5
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = -3.333333333333
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 101
$ws.Range("K20").Value = 4.950495049504
$ws.Range("L20").Value = 7.070707070707
$ws.Range("M20").Value = 17.777777777777
$ws.Range("N20").Value = -92.205882352941
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -18.75
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -17.796610169491
$ws.Range("I21").Value = 426
$ws.Range("J21").Value = 486
$ws.Range("K21").Value = -12.345679012345
$ws.Range("L21").Value = -11.983471074380
$ws.Range("M21").Value = -13.238289205702
$ws.Range("N21").Value = -83.857521788556
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -40
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -3.125
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = -6.666666666666
$ws.Range("I24").Value = 494
$ws.Range("J24").Value = 552
$ws.Range("K24").Value = -10.507246376811
$ws.Range("L24").Value = -10.990990990991
$ws.Range("M24").Value = 69.178082191780
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 73
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 43.137254901960
$ws.Range("I25").Value = 254
$ws.Range("J25").Value = 258
$ws.Range("K25").Value = -1.550387596899
$ws.Range("L25").Value = 52.095808383233
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 57
$ws.Range("H26").Value = -10.526315789473
$ws.Range("I26").Value = 198
$ws.Range("J26").Value = 237
$ws.Range("K26").Value = -16.455696202531
$ws.Range("L26").Value = -1.492537313432
$ws.Range("M26").Value = -9.589041095890
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -7.692307692307
$ws.Range("L27").Value = -14.285714285714
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 60
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 45.454545454545
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 0
